# Daily price-data update: insert a new "today" row at the top of the
# table (row 2, just below the header) and push the existing history
# down by one row. The new row repeats the same price figures
# (783.5 / 1112 / 3610) that every other row already uses, dated one
# day after the previous most-recent entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 2:end down by one to make room for the new latest entry.
$ws.Rows("2:2").Insert()

# Write the new date as literal text (not an auto-converted date
# serial) so it matches the existing column formatting, then drop the
# quote-prefix style Excel adds for text-that-looks-like-a-date so the
# cell keeps the sheet's default (unstyled) look.
$ws.Range("A2").Value = "'2026-01-25"
$ws.Range("A2").Style = "Normal"

$ws.Range("B2").Value = 783.5
$ws.Range("C2").Value = 1112
$ws.Range("D2").Value = 3610
